$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row above row 11 (current "Description" row) to make room for "Jurisdiction".
$ws.Rows.Item(11).Insert()

# The inserted row picks up a generic style; copy the format from the row below
# (the old "Description" row, now pushed down to row 12) so it matches the rest
# of the table (border + wrap formatting).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Populate the new "Jurisdiction" property row (value left blank).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Update the Version and Date property values.
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"
